$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.113.92'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '2.281.27'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.45'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.60'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.620'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.99%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.85'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0899'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.22'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.87%  '
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.951'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.04'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("D16").Value = '2.626.75'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = '2.280.64'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").Value = '42.093.09'
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.29'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000105'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.67'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +28.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.60'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '263.55'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.21'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.73%  '
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.75'
$ws.Range("D27").ClearFormats()
$ws.Range("E28").Value = '  +1.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.38'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.34'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '163.65'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.04'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0865'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.30%  '
$ws.Range("E34").Value = '  +1.53%  '
$ws.Range("E35").Value = '  -3.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.51'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -13.81%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.54'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0353'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.53%  '
$ws.Range("E39").Value = '  -5.17%  '
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.52'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.28%  '
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '68.13'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.17%  '
$ws.Range("E44").Value = '  -1.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.20'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -7.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '115.02'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.90'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '78.77'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.92'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.01%  '
$ws.Range("D50").Value = '1.607.54'
$ws.Range("E50").Value = '  +4.24%  '
$ws.Range("E51").Value = '  -2.08%  '
